# LOQ4247 worksheet update
#
# The course-description sheet had its "Docentes responsaveis" value row
# (old row 13, holding only the teacher name with no label in column A)
# removed, which shifts every row below it up by one. A handful of cells
# that land on the newly-shifted rows pick up different (re-keyed) values
# rather than simply inheriting the row that used to be one below them, so
# those are patched explicitly afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old, label-less row 13 ("8767640 - Eduardo Ferro dos Santos"
# with nothing in column A). Everything below shifts up by one row.
$ws.Rows.Item(13).Delete()

# Row 10 ("Objetivos:") now shows the teacher name instead of the old
# Portuguese objectives paragraph.
$ws.Range("B10").Value = "8767640 - Eduardo Ferro dos Santos"
$ws.Range("C10").Value = "8767640 - Eduardo Ferro dos Santos"

# Row 13 ("Programa resumido:") now shows "Semestral" instead of the long
# Portuguese summary paragraph.
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 ("Programa:") now shows the activation date instead of the long
# Portuguese program paragraph. Entered via copy/paste-values from the
# existing "01/01/2018" text cell (B8/C8) so it lands back in the shared
# string table as plain text instead of Excel auto-converting the literal
# "01/01/2018" into a date serial number.
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = 0

# Row 18 ("Metodo:") now shows the teacher name instead of the
# "Aulas expositivas e praticas." text.
$ws.Range("B18").Value = "8767640 - Eduardo Ferro dos Santos"
$ws.Range("C18").Value = "8767640 - Eduardo Ferro dos Santos"

# Row 19 ("Criterio:") now shows "Aulas expositivas e praticas." instead of
# the evaluation-criteria paragraph.
$ws.Range("B19").Value = "Aulas expositivas e práticas."
$ws.Range("C19").Value = "Aulas expositivas e práticas."

# Row 20 ("Norma de recuperacao:") now shows the evaluation-criteria
# paragraph instead of the recovery-norm paragraph.
$ws.Range("B20").Value = "Exercícios de aprendizado e exercícios de avaliação farão parte da composição de notas individuais (NI), com aplicação de trabalhos práticos em grupo (NG). Sendo: Nota Final = (NI+NG)/2"
$ws.Range("C20").Value = "Exercícios de aprendizado e exercícios de avaliação farão parte da composição de notas individuais (NI), com aplicação de trabalhos práticos em grupo (NG). Sendo: Nota Final = (NI+NG)/2"

# Row 21 ("Bibliografia:") now shows the recovery-norm paragraph instead of
# the full bibliography text (which is dropped entirely).
$ws.Range("B21").Value = "A recuperação deverá consistir de uma prova englobando a matéria toda do semestre. - A média final (pós-recuperação) deverá ser composta por uma média simples entre a nota do semestre (nota final) e a da prova de recuperação."
$ws.Range("C21").Value = "A recuperação deverá consistir de uma prova englobando a matéria toda do semestre. - A média final (pós-recuperação) deverá ser composta por uma média simples entre a nota do semestre (nota final) e a da prova de recuperação."
